$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price-report row per (market, product) observation,
# ordered (mostly) by date. Two new weekly rows were added:
#   - one inserted at row 5 (pushing the old rows 5-10 down to 6-11)
#   - one inserted at row 9 of the new layout (pushing what is now at 9-11 down to 10-12)
# giving a final block of 8 data rows (5-12) instead of the original 6 (5-10).

$ws.Rows("5:5").Insert()
$ws.Rows("9:9").Insert()

# New row inserted at position 5 (newest weekly observation).
$row5 = @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44965, 4, 100112039, "Ciboulette", "Sin especificar", "Primera", 1120, 2000, 2500, 2250, "`$/docena de atados", "Provincia del Elquí", 750, 3, "Hortaliza")
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# New row inserted at position 9 (second newest weekly observation).
$row9 = @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44964, 4, 100112039, "Ciboulette", "Sin especificar", "Primera", 1000, 2000, 2500, 2250, "`$/docena de atados", "Provincia del Elquí", 750, 3, "Hortaliza")
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}
